# Update market-price-derived profit figures per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1607
$ws.Range("I11").Value = 1607
$ws.Range("K11").Value = 1607
$ws.Range("M11").Value = -1467
$ws.Range("H12").Value = 627
$ws.Range("I12").Value = 439.5
$ws.Range("K12").Value = 439.5
$ws.Range("M12").Value = -269.5
$ws.Range("H28").Value = 1774.2222
$ws.Range("I28").Value = 1213.3846
$ws.Range("K28").Value = 1213.3846
$ws.Range("M28").Value = -728.3846000000001
$ws.Range("H86").Value = 65617616
$ws.Range("I86").Value = 125003200
$ws.Range("K86").Value = 125003200
$ws.Range("M86").Value = -125002077
$ws.Range("H89").Value = 65617616
$ws.Range("I89").Value = 125003200
$ws.Range("K89").Value = 625016000
$ws.Range("M89").Value = -625010384
$ws.Range("H103").Value = 1094.9445
$ws.Range("J103").Value = 1168
$ws.Range("L103").Value = 3504
$ws.Range("N103").Value = -4676
$ws.Range("H106").Value = 1664.375
$ws.Range("I106").Value = 1254.1111
$ws.Range("K106").Value = 1254.1111
$ws.Range("M106").Value = -623.1111000000001
$ws.Range("H112").Value = 4278.1665
$ws.Range("J112").Value = 4278.1665
$ws.Range("L112").Value = 12834.4995
$ws.Range("N112").Value = -15050.4995
$ws.Range("H113").Value = 25013452
$ws.Range("I113").Value = 2484.6667
$ws.Range("K113").Value = 2484.6667
$ws.Range("M113").Value = 769.3332999999998
$ws.Range("H116").Value = 14714727
$ws.Range("J116").Value = 10431.909
$ws.Range("L116").Value = 10431.909
$ws.Range("N116").Value = -17315.909
$ws.Range("H132").Value = 2132.2444
$ws.Range("I132").Value = 2207.1428
$ws.Range("K132").Value = 6621.428400000001
$ws.Range("M132").Value = -4091.428400000001
$ws.Range("H137").Value = 3442.3901
$ws.Range("I137").Value = 3013.7222
$ws.Range("J137").Value = 3777.8696
$ws.Range("K137").Value = 9041.1666
$ws.Range("L137").Value = 11333.6088
$ws.Range("M137").Value = -6491.1666
$ws.Range("N137").Value = -16433.6088
$ws.Range("H141").Value = 1897.3143
$ws.Range("I141").Value = 1804.8148
$ws.Range("J141").Value = 2209.5
$ws.Range("K141").Value = 5414.4444
$ws.Range("L141").Value = 6628.5
$ws.Range("M141").Value = -234.4444000000003
$ws.Range("N141").Value = -16988.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2022084.6
$ws.Range("I32").Value = 2088237.1
$ws.Range("K32").Value = 2088237.1
$ws.Range("M32").Value = -2087950.1
$ws.Range("H46").Value = 4581.375
$ws.Range("J46").Value = 4664.4287
$ws.Range("L46").Value = 4664.4287
$ws.Range("N46").Value = -5302.4287
$ws.Range("H74").Value = 19398.12
$ws.Range("I74").Value = 24611.092
$ws.Range("K74").Value = 24611.092
$ws.Range("M74").Value = -23737.092
$ws.Range("H77").Value = 19398.12
$ws.Range("I77").Value = 24611.092
$ws.Range("K77").Value = 123055.46
$ws.Range("M77").Value = -118687.46
$ws.Range("H122").Value = 4110.8184
$ws.Range("I122").Value = 2902.7827
$ws.Range("J122").Value = 6889.3
$ws.Range("K122").Value = 8708.348100000001
$ws.Range("L122").Value = 20667.9
$ws.Range("M122").Value = -6258.348100000001
$ws.Range("N122").Value = -25567.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7578043.5
$ws.Range("I20").Value = 11906947
$ws.Range("J20").Value = 2462.375
$ws.Range("K20").Value = 11906947
$ws.Range("L20").Value = 2462.375
$ws.Range("M20").Value = -11906700
$ws.Range("N20").Value = -2956.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7167.45
$ws.Range("I31").Value = 4004.525
$ws.Range("K31").Value = 4004.525
$ws.Range("M31").Value = -3709.525
$ws.Range("H34").Value = 7167.45
$ws.Range("I34").Value = 4004.525
$ws.Range("K34").Value = 4004.525
$ws.Range("M34").Value = -3802.525
$ws.Range("H134").Value = 3210.6
$ws.Range("I134").Value = 1233.9403
$ws.Range("K134").Value = 3701.8209
$ws.Range("M134").Value = -1166.8209
$ws.Range("H139").Value = 93750
$ws.Range("J139").Value = 93750
$ws.Range("L139").Value = 93750
$ws.Range("N139").Value = -104030

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2867.9546
$ws.Range("I5").Value = 1002.875
$ws.Range("J5").Value = 3933.7144
$ws.Range("K5").Value = 3008.625
$ws.Range("L5").Value = 11801.1432
$ws.Range("M5").Value = -2896.625
$ws.Range("N5").Value = -12025.1432
$ws.Range("H76").Value = 2999
$ws.Range("I76").Value = 2999
$ws.Range("K76").Value = 8997
$ws.Range("M76").Value = -8614
$ws.Range("H79").Value = 2999
$ws.Range("I79").Value = 2999
$ws.Range("K79").Value = 8997
$ws.Range("M79").Value = -7671
$ws.Range("H122").Value = 1572703.2
$ws.Range("I122").Value = 3143640.5
$ws.Range("J122").Value = 1765.8889
$ws.Range("K122").Value = 28292764.5
$ws.Range("L122").Value = 15893.0001
$ws.Range("M122").Value = -28290314.5
$ws.Range("N122").Value = -20793.0001
$ws.Range("H132").Value = 4713.6665
$ws.Range("I132").Value = 3094.65
$ws.Range("J132").Value = 6185.5
$ws.Range("K132").Value = 27851.85
$ws.Range("L132").Value = 55669.5
$ws.Range("M132").Value = -25321.85
$ws.Range("N132").Value = -60729.5
$ws.Range("H135").Value = 2867.9546
$ws.Range("I135").Value = 1002.875
$ws.Range("J135").Value = 3933.7144
$ws.Range("K135").Value = 9025.875
$ws.Range("L135").Value = 35403.4296
$ws.Range("M135").Value = -6490.875
$ws.Range("N135").Value = -40473.4296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 79999.75
$ws.Range("J52").Value = 79999.75
$ws.Range("L52").Value = 79999.75
$ws.Range("N52").Value = -80517.75
$ws.Range("H97").Value = 673.15
$ws.Range("I97").Value = 524.8333
$ws.Range("K97").Value = 524.8333
$ws.Range("M97").Value = -28.83330000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4275.6855
$ws.Range("I122").Value = 2767.4783
$ws.Range("J122").Value = 7166.4165
$ws.Range("K122").Value = 8302.4349
$ws.Range("L122").Value = 21499.2495
$ws.Range("M122").Value = -5852.4349
$ws.Range("N122").Value = -26399.2495
$ws.Range("H132").Value = 5751305
$ws.Range("I132").Value = 8476589
$ws.Range("K132").Value = 25429767
$ws.Range("M132").Value = -25427237

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11796090
$ws.Range("I81").Value = 46858.91
$ws.Range("J81").Value = 33336348
$ws.Range("K81").Value = 93717.82000000001
$ws.Range("L81").Value = 66672696
$ws.Range("M81").Value = -92656.82000000001
$ws.Range("N81").Value = -66674818
$ws.Range("H84").Value = 11796090
$ws.Range("I84").Value = 46858.91
$ws.Range("J84").Value = 33336348
$ws.Range("K84").Value = 468589.1
$ws.Range("L84").Value = 333363480
$ws.Range("M84").Value = -463285.1
$ws.Range("N84").Value = -333374088
$ws.Range("H107").Value = 9009792
$ws.Range("I107").Value = 434.72223
$ws.Range("J107").Value = 17544972
$ws.Range("K107").Value = 1304.16669
$ws.Range("L107").Value = 52634916
$ws.Range("M107").Value = 615.83331
$ws.Range("N107").Value = -52638756
